$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update tariff amounts: D5 5 -> 10, D6 10 -> 20
$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 20

# Update the active selection to D6
$ws.Range("D6").Select()
